$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 6
$ws.Range("E2").Value = "Git"
$ws.Range("F2").Value = "GitKraken"
$ws.Range("G2").Value = "GitKraken"

$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = "Git"
$ws.Range("F3").Value = "GitUp"
$ws.Range("G3").Value = "GitUp"

$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = "GitKraken"
$ws.Range("F4").Value = "GitUp"
$ws.Range("G4").Value = "GitUp"

$ws.Range("D22").Select()
